# Fill in "My answer" (column C) on sheet "2" ("Exam Attempt2") to match the
# grader state: most answers reuse the existing question-bank text, a handful
# are near-miss typos that keep D (the IF(B=C,...) check) at FALSE.
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value = 'The hacker still has the opportunity to connect to the network after sniffing the SSID from a successful wireless association'
$ws2.Range("C3").Value = 'White hat'
$ws2.Range("C4").Value = 'Logic tier'
$ws2.Range("C5").Value = 'Jacob inadvertently provided the answers to his security questions when responding to Jane''s post'
$ws2.Range("C6").Value = 'Directory traversal'
$ws2.Range("C7").Value = 'Initial intrusion'
$ws2.Range("C8").Value = 'UDP hijacking'
$ws2.Range("C9").Value = 'Weaponization'
$ws2.Range("C10").Value = 'php.ini'
$ws2.Range("C12").Value = 'Use your colleague''s public key to encrypt the message.'
$ws2.Range("C13").Value = 'Social engineering'
$ws2.Range("C14").Value = 'Jennys will sign the message with her private key, and Molly will verify that the message came from Jennys by using Jenny’s public key'
$ws2.Range("C15").Value = 'DNS tunnelling '
$ws2.Range("C11").Value = 'hping2 -1 target.domain.com '
$ws2.Range("C16").Value = 'A browser makes a request to a server without the user''s knowledge'
$ws2.Range("C17").Value = 'Agent-based scanner'

# Row 17 ends up referencing a deleted cell (#REF!); Excel then re-anchors the
# next few members of that shared-formula block one row back, so rows 18-20
# stay FALSE even though their "My answer" matches the correct answer.
$ws2.Range("D17").Formula = "=IF(B17=#REF!,TRUE,FALSE)"
$ws2.Range("C18").Value = 'Wireshark with Airpcap'
$ws2.Range("D18").Formula = "=IF(B18=C17,TRUE,FALSE)"
$ws2.Range("C19").Value = 'ARP ping scan'
$ws2.Range("D19").Formula = "=IF(B19=C18,TRUE,FALSE)"
$ws2.Range("C20").Value = '.bash_history '
$ws2.Range("D20").Formula = "=IF(B20=C19,TRUE,FALSE)"

$ws2.Range("C21").Value = 'Whois footprinting '
$ws2.Range("C22").Value = 'Remote-access policy '
$ws2.Range("C23").Value = 'Impersonation attack'
$ws2.Range("C24").Value = 'Server Message Block (SMB)'
$ws2.Range("C25").Value = 'He is scanning from 192.168.1.64 to 192.168.1.78 because of the mask /28 and the servers are not in that range'
$ws2.Range("C26").Value = 'Docker daemon'
$ws2.Range("C27").Value = 'Burp suite'
$ws2.Range("C28").Value = 'n'
$ws2.Range("C29").Value = 'Bug bounty program '
$ws2.Range("C30").Value = 'Quid pro quo'
$ws2.Range("C31").Value = 'Adware'
$ws2.Range("C32").Value = 128
$ws2.Range("C33").Value = 'The hacker makes a request to the DNS resolver'
$ws2.Range("C34").Value = 'Code Emulation'
$ws2.Range("C35").Value = 'Advanced persistent threat'
$ws2.Range("C36").Value = 'Website mirroring'
$ws2.Range("C37").Value = 'It''s a stateful firewall'
$ws2.Range("C38").Value = 'Actions on objectives'
$ws2.Range("C39").Value = 'Agent Smith attack'
$ws2.Range("C40").Value = 'iOS trustjacking'
$ws2.Range("C41").Value = 'nmap -sn -PS < target IP address > '
$ws2.Range("C42").Value = '10.1.155.200'
$ws2.Range("C43").Value = 'Fileless malware '
$ws2.Range("C44").Value = 'STP attack'
$ws2.Range("C45").Value = 'Untethered jailbreaking'
$ws2.Range("C46").Value = 'Gaining access'
$ws2.Range("C47").Value = 'Webroot'
$ws2.Range("C48").Value = 'FCC ID search'
$ws2.Range("C49").Value = 'Union-based and error-based'
$ws2.Range("C50").Value = 'Bluesmacking'
$ws2.Range("C51").Value = 'Private keys'
$ws2.Range("C52").Value = 'VRFY'
$ws2.Range("C53").Value = 'JXplorer'
$ws2.Range("C54").Value = 'Host-based assessment'
$ws2.Range("C55").Value = 'Cloudborne attack'
$ws2.Range("C56").Value = 'Server-side request forgery attack'
$ws2.Range("C57").Value = 'Topological scanning technique'
$ws2.Range("C58").Value = 'Evil twin attack '
$ws2.Range("C59").Value = 'WPA3-Enterprise'
$ws2.Range("C60").Value = 'Phishing attack: an attacker provides the victim with a URL that is either misspelled or looks similar to the legitimate website''s domain name. Pharming attack: a victim is redirected to a fake website by modifying their host configuration file or exploiting DNS vulnerabilities'
$ws2.Range("C61").Value = 53
$ws2.Range("C62").Value = 'Dragonblood'

# Switch the active sheet/selection to sheet 2, mirroring the saved view state
$ws2.Activate()
$ws2.Range("C63").Select()
